$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Rename the individual id used for the "TestScenario" population rows
# (B4, B5) from "Indiv" to "Indiv1", matching the value already used in B2.
$ws.Range("B4").Value = "Indiv1"
$ws.Range("B5").Value = "Indiv1"

# Update the last active selection on the Scenarios sheet.
$ws.Range("E19").Select() | Out-Null
